# AFA 2020.xlsx — Contest 24 KXI vs KKR and Contest 25 CSK vs RCB
# - fills in the scores for rows 33 (Contest 24) and 34 (Contest 25)
# - adds two new contest rows (33 "RR vs RCB" and 34 "DC vs CSK") by
#   inserting two fresh template rows right before row 44 (pushing the
#   "Total" block down and widening the SUM ranges the same way Excel does)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Contest 24 (row 33, "KXI vs KKR") - enter raw scores.
#    E33/T33 tie at 40, so D33 and S33 are given a manually split score
#    (-17.5, halfway between the rank-4 and rank-5 payouts) instead of
#    the RANK/VLOOKUP formula, matching what a human would enter to
#    break the tie fairly.
# ---------------------------------------------------------------------
$ws.Range("E33").Value = 40
$ws.Range("H33").Value = 0
$ws.Range("K33").Value = 60
$ws.Range("N33").Value = 80
$ws.Range("Q33").Value = 100
$ws.Range("T33").Value = 40
$ws.Range("D33").Value = -17.5
$ws.Range("S33").Value = -17.5

# ---------------------------------------------------------------------
# 2) Contest 25 (row 34, "CSK vs RCB") - enter raw scores; no ties here
#    so every formula cell (D/G/J/M/P/S) recomputes on its own.
# ---------------------------------------------------------------------
$ws.Range("E34").Value = 60
$ws.Range("H34").Value = 40
$ws.Range("K34").Value = 100
$ws.Range("N34").Value = 80
$ws.Range("Q34").Value = 20
$ws.Range("T34").Value = 0

# ---------------------------------------------------------------------
# 3) Make room for two new contests by inserting two rows right at the
#    existing blank template rows 42:43 (not at 44) -- this is the spot
#    where Excel widens the SUM(...10:...42) ranges below to .../...44
#    because the insertion happens inside the summed range.
# ---------------------------------------------------------------------
$ws.Rows("42:43").Insert()

$cols = @("A","B","C","D","E","G","H","J","K","M","N","P","Q","S","T")
foreach ($col in $cols) {
  $ws.Range($col + "41").Copy()
  $ws.Range($col + "42").PasteSpecial(-4122)
  $ws.Range($col + "43").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

foreach ($r in 42,43) {
  $ws.Range("D" + $r).Formula = "=IF(ISERROR(VLOOKUP(RANK(E" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(E" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE))"
  $ws.Range("G" + $r).Formula = "=IF(ISERROR(VLOOKUP(RANK(H" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(H" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE))"
  $ws.Range("J" + $r).Formula = "=IF(ISERROR(VLOOKUP(RANK(K" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(K" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE))"
  $ws.Range("M" + $r).Formula = "=IF(ISERROR(VLOOKUP(RANK(N" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(N" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE))"
  $ws.Range("P" + $r).Formula = "=IF(ISERROR(VLOOKUP(RANK(Q" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(Q" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE))"
  $ws.Range("S" + $r).Formula = "=IF(ISERROR(VLOOKUP(RANK(T" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(T" + $r + ", (`$T" + $r + ",`$Q" + $r + ",`$N" + $r + ",`$K" + $r + ",`$H" + $r + ",`$E" + $r + "), 0),  score, 2, FALSE))"
}

# ---------------------------------------------------------------------
# 4) New contest 33 "RR vs RCB" and contest 34 "DC vs CSK" - no scores
#    yet, so only the match-number / format / name columns are filled.
# ---------------------------------------------------------------------
$ws.Range("A42").Value = 33
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = "RR vs RCB"

$ws.Range("A43").Value = 34
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = "DC vs CSK"
